$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe is Excel's "force text" qualifier: it is stripped
# from the cell's stored value but prevents numeric-looking text (e.g.
# "211.76", "0.100", "27.953.49") from being auto-coerced into a real
# number, which would corrupt the original formatting / drop trailing
# zeros. Every written value below keeps the exact source text.
$apos = [string][char]39

$ws.Range('D2').Value = $apos + '27.953.49'
$ws.Range('E2').Value = $apos + '  +0.16%  '
$ws.Range('D3').Value = $apos + '1.631.47'
$ws.Range('E3').Value = $apos + '  -0.48%  '
$ws.Range('E4').Value = $apos + '  -0.03%  '
$ws.Range('D5').Value = $apos + '211.76'
$ws.Range('E5').Value = $apos + '  -0.90%  '
$ws.Range('E6').Value = $apos + '  -0.14%  '
$ws.Range('E7').Value = $apos + '  -0.06%  '
$ws.Range('D8').Value = $apos + '23.43'
$ws.Range('E8').Value = $apos + '  -1.01%  '
$ws.Range('E9').Value = $apos + '  -1.79%  '
$ws.Range('E10').Value = $apos + '  -0.22%  '
$ws.Range('D11').Value = $apos + '0.0881'
$ws.Range('E11').Value = $apos + '  +0.65%  '
$ws.Range('D12').Value = $apos + '1.863.29'
$ws.Range('E12').Value = $apos + '  -0.46%  '
$ws.Range('D13').Value = $apos + '1.626.46'
$ws.Range('E13').Value = $apos + '  -0.94%  '
$ws.Range('E14').Value = $apos + '  -1.26%  '
$ws.Range('E15').Value = $apos + '  -1.91%  '
$ws.Range('E16').Value = $apos + '  -0.93%  '
$ws.Range('D17').Value = $apos + '27.949.04'
$ws.Range('D18').Value = $apos + '231.00'
$ws.Range('E18').Value = $apos + '  -0.48%  '
$ws.Range('D19').Value = $apos + '0.0₃0725'
$ws.Range('E19').Value = $apos + '  +0.07%  '
$ws.Range('D20').Value = $apos + '7.64'
$ws.Range('E20').Value = $apos + '  +0.44%  '
$ws.Range('D21').Value = $apos + '0.999'
$ws.Range('E21').Value = $apos + '  -0.08%  '
$ws.Range('D22').Value = $apos + '10.37'
$ws.Range('E22').Value = $apos + '  -5.31%  '
$ws.Range('E24').Value = $apos + '  -1.57%  '
$ws.Range('D25').Value = $apos + '155.00'
$ws.Range('E25').Value = $apos + '  +2.13%  '
$ws.Range('D26').Value = $apos + '6.94'
$ws.Range('E26').Value = $apos + '  +0.28%  '
$ws.Range('E27').Value = $apos + '  -0.07%  '
$ws.Range('D28').Value = $apos + '15.58'
$ws.Range('E28').Value = $apos + '  -0.79%  '
$ws.Range('E29').Value = $apos + '  -0.12%  '
$ws.Range('E30').Value = $apos + '  -0.57%  '
$ws.Range('E31').Value = $apos + '  -0.56%  '
$ws.Range('E32').Value = $apos + '  +1.83%  '
$ws.Range('D33').Value = $apos + '1.401.53'
$ws.Range('E33').Value = $apos + '  -1.19%  '
$ws.Range('E34').Value = $apos + '  -1.67%  '
$ws.Range('D35').Value = $apos + '1.58'
$ws.Range('E35').Value = $apos + '  +0.35%  '
$ws.Range('E36').Value = $apos + '  +10.90%  '
$ws.Range('E37').Value = $apos + '  +0.54%  '
$ws.Range('D38').Value = $apos + '0.0171'
$ws.Range('E38').Value = $apos + '  +2.45%  '
$ws.Range('E39').Value = $apos + '  +0.03%  '
$ws.Range('D40').Value = $apos + '0.867'
$ws.Range('E40').Value = $apos + '  -2.93%  '
$ws.Range('D41').Value = $apos + '1.03'
$ws.Range('E41').Value = $apos + '  -0.01%  '
$ws.Range('E42').Value = $apos + '  -0.08%  '
$ws.Range('D43').Value = $apos + '66.52'
$ws.Range('E43').Value = $apos + '  -0.86%  '
$ws.Range('D44').Value = $apos + '1.84'
$ws.Range('E44').Value = $apos + '  +0.79%  '
$ws.Range('D45').Value = $apos + '5.46'
$ws.Range('E45').Value = $apos + '  +0.38%  '
$ws.Range('E46').Value = $apos + '  -0.35%  '
$ws.Range('D47').Value = $apos + '1.773.11'
$ws.Range('D48').Value = $apos + '88.30'
$ws.Range('E48').Value = $apos + '  -0.05%  '
$ws.Range('B49').Value = $apos + 'Algorand'
$ws.Range('C49').Value = $apos + 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = $apos + '0.100'
$ws.Range('E49').Value = $apos + '  -0.18%  '
$ws.Range('B50').Value = $apos + 'Cronos'
$ws.Range('C50').Value = $apos + 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = $apos + '0.0505'
$ws.Range('E50').Value = $apos + '  -0.29%  '
$ws.Range('B51').Value = $apos + 'EnergySwap'
$ws.Range('C51').Value = $apos + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = $apos + '7.54'
$ws.Range('E51').Value = $apos + '  -1.24%  '
